# Add a new "Save" column (H) to the s_vals sheet, mirroring the existing
# "Win" (F) style column: a bold/bordered header in row 1 plus a 0/1 value
# for each data row (2-17).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give the new header cell the same formatting as the other header cells
# (bold font + border + centered alignment) by copying G1's format onto H1,
# then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("H1").Value = "Save"

# Save indicator values for rows 2-17 (one per data row, in sheet order).
$saveValues = @(1, 0, 1, 0, 0, 1, 0, 0, 0, 0, 1, 0, 0, 1, 0, 0)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
